$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Reserved / bought" column (E) with "Y" for the reserved/bought rows.
$ws.Range("E2").Value = "Y"
$ws.Range("E3").Value = "Y"
$ws.Range("E9").Value = "Y"

# Update the selected cell to match the author's last selection.
$ws.Range("E3").Select()
